$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.026992
$ws.Range("H2").Value = 0.08097599999999999
$ws.Range("I2").Value = 0.004182906599909731
$ws.Range("J2").Value = 0.00420788870005516
$ws.Range("M2").Value = 89.38217433333334
$ws.Range("N2").Value = 268.146523
$ws.Range("O2").Value = 0.2143552015363441
$ws.Range("P2").Value = 0.2175965347165783
$ws.Range("Q2").Value = 2.412603649605333
$ws.Range("R2").Value = 21.713432846448
$ws.Range("S2").Value = 0.0008966277872313542
$ws.Range("T2").Value = 0.0009156219996050504
$ws.Range("G3").Value = 0.026992
$ws.Range("H3").Value = 0.08097599999999999
$ws.Range("I3").Value = 0.004182906599909731
$ws.Range("J3").Value = 0.00420788870005516
$ws.Range("O3").Value = 0.2934277926151677
$ws.Range("P3").Value = 0.2978648075949286
$ws.Range("Q3").Value = 3.302578889082667
$ws.Range("R3").Value = 29.723210001744
$ws.Range("S3").Value = 0.001227381050326928
$ws.Range("T3").Value = 0.001253381958022804
$ws.Range("G4").Value = 0.026992
$ws.Range("H4").Value = 0.08097599999999999
$ws.Range("I4").Value = 0.004182906599909731
$ws.Range("J4").Value = 0.00420788870005516
$ws.Range("M4").Value = 90.33462533333334
$ws.Range("N4").Value = 271.003876
$ws.Range("O4").Value = 0.2166393574945233
$ws.Range("P4").Value = 0.2199152301234996
$ws.Range("Q4").Value = 2.438312206997333
$ws.Range("R4").Value = 21.944809862976
$ws.Range("S4").Value = 0.0009061821982640449
$ws.Range("T4").Value = 0.0009253788118067042
$ws.Range("G5").Value = 0.026992
$ws.Range("H5").Value = 0.08097599999999999
$ws.Range("I5").Value = 0.004182906599909731
$ws.Range("J5").Value = 0.00420788870005516
$ws.Range("M5").Value = 18.634161
$ws.Range("N5").Value = 37.268322
$ws.Range("O5").Value = 0.0446882095496985
$ws.Range("P5").Value = 0.03024263611988591
$ws.Range("Q5").Value = 0.502973273712
$ws.Range("R5").Value = 3.017839642272
$ws.Range("S5").Value = 0.0001869266066635829
$ws.Range("T5").Value = 0.000127257646788748
$ws.Range("G6").Value = 0.026992
$ws.Range("H6").Value = 0.08097599999999999
$ws.Range("I6").Value = 0.004182906599909731
$ws.Range("J6").Value = 0.00420788870005516
$ws.Range("M6").Value = 96.27664699999998
$ws.Range("N6").Value = 288.829941
$ws.Range("O6").Value = 0.2308894388042666
$ws.Range("P6").Value = 0.2343807914451077
$ws.Range("Q6").Value = 2.598699255823999
$ws.Range("R6").Value = 23.388293302416
$ws.Range("S6").Value = 0.0009657889574238204
$ws.Range("T6").Value = 0.0009862482838318538
$ws.Range("I7").Value = 0.976381346197431
$ws.Range("J7").Value = 0.9822127115383066
$ws.Range("M7").Value = 89.38217433333334
$ws.Range("N7").Value = 268.146523
$ws.Range("O7").Value = 0.2143552015363441
$ws.Range("P7").Value = 0.2175965347165783
$ws.Range("Q7").Value = 563.1541472365952
$ws.Range("R7").Value = 5068.387325129357
$ws.Range("S7").Value = 0.2092924202404773
$ws.Range("T7").Value = 0.2137260823853097
$ws.Range("I8").Value = 0.976381346197431
$ws.Range("J8").Value = 0.9822127115383066
$ws.Range("O8").Value = 0.2934277926151677
$ws.Range("P8").Value = 0.2978648075949286
$ws.Range("S8").Value = 0.286497423165338
$ws.Range("T8").Value = 0.2925666003396508
$ws.Range("I9").Value = 0.976381346197431
$ws.Range("J9").Value = 0.9822127115383066
$ws.Range("M9").Value = 90.33462533333334
$ws.Range("N9").Value = 271.003876
$ws.Range("O9").Value = 0.2166393574945233
$ws.Range("P9").Value = 0.2199152301234996
$ws.Range("Q9").Value = 569.1550834936315
$ws.Range("R9").Value = 5122.395751442684
$ws.Range("S9").Value = 0.2115226275098491
$ws.Range("T9").Value = 0.2160035344881732
$ws.Range("I10").Value = 0.976381346197431
$ws.Range("J10").Value = 0.9822127115383066
$ws.Range("M10").Value = 18.634161
$ws.Range("N10").Value = 37.268322
$ws.Range("O10").Value = 0.0446882095496985
$ws.Range("P10").Value = 0.03024263611988591
$ws.Range("Q10").Value = 117.404897852333
$ws.Range("R10").Value = 704.4293871139979
$ws.Range("S10").Value = 0.04363273419928752
$ws.Range("T10").Value = 0.02970470162737947
$ws.Range("I11").Value = 0.976381346197431
$ws.Range("J11").Value = 0.9822127115383066
$ws.Range("M11").Value = 96.27664699999998
$ws.Range("N11").Value = 288.829941
$ws.Range("O11").Value = 0.2308894388042666
$ws.Range("P11").Value = 0.2343807914451077
$ws.Range("Q11").Value = 606.5929078642242
$ws.Range("R11").Value = 5459.336170778018
$ws.Range("S11").Value = 0.2254361410824792
$ws.Range("T11").Value = 0.2302117926977935
$ws.Range("G12").Value = 0.1149325
$ws.Range("H12").Value = 0.229865
$ws.Range("I12").Value = 0.01781090370458377
$ws.Range("J12").Value = 0.01194485200600399
$ws.Range("M12").Value = 89.38217433333334
$ws.Range("N12").Value = 268.146523
$ws.Range("O12").Value = 0.2143552015363441
$ws.Range("P12").Value = 0.2175965347165783
$ws.Range("Q12").Value = 10.27291675156583
$ws.Range("R12").Value = 61.637500509395
$ws.Range("S12").Value = 0.003817859853140472
$ws.Range("T12").Value = 0.002599158404208838
$ws.Range("G13").Value = 0.1149325
$ws.Range("H13").Value = 0.229865
$ws.Range("I13").Value = 0.01781090370458377
$ws.Range("J13").Value = 0.01194485200600399
$ws.Range("O13").Value = 0.2934277926151677
$ws.Range("P13").Value = 0.2978648075949286
$ws.Range("Q13").Value = 14.06244991736417
$ws.Range("R13").Value = 84.37469950418499
$ws.Range("S13").Value = 0.005226214158517328
$ws.Range("T13").Value = 0.003557951044518276
$ws.Range("G14").Value = 0.1149325
$ws.Range("H14").Value = 0.229865
$ws.Range("I14").Value = 0.01781090370458377
$ws.Range("J14").Value = 0.01194485200600399
$ws.Range("M14").Value = 90.33462533333334
$ws.Range("N14").Value = 271.003876
$ws.Range("O14").Value = 0.2166393574945233
$ws.Range("P14").Value = 0.2199152301234996
$ws.Range("Q14").Value = 10.38238432612333
$ws.Range("R14").Value = 62.29430595674
$ws.Range("S14").Value = 0.003858542734957852
$ws.Range("T14").Value = 0.002626854877691514
$ws.Range("G15").Value = 0.1149325
$ws.Range("H15").Value = 0.229865
$ws.Range("I15").Value = 0.01781090370458377
$ws.Range("J15").Value = 0.01194485200600399
$ws.Range("M15").Value = 18.634161
$ws.Range("N15").Value = 37.268322
$ws.Range("O15").Value = 0.0446882095496985
$ws.Range("P15").Value = 0.03024263611988591
$ws.Range("Q15").Value = 2.1416707091325
$ws.Range("R15").Value = 8.566682836529999
$ws.Range("S15").Value = 0.0007959373970199408
$ws.Range("T15").Value = 0.0003612438127234681
$ws.Range("G16").Value = 0.1149325
$ws.Range("H16").Value = 0.229865
$ws.Range("I16").Value = 0.01781090370458377
$ws.Range("J16").Value = 0.01194485200600399
$ws.Range("M16").Value = 96.27664699999998
$ws.Range("N16").Value = 288.829941
$ws.Range("O16").Value = 0.2308894388042666
$ws.Range("P16").Value = 0.2343807914451077
$ws.Range("Q16").Value = 11.0653157313275
$ws.Range("R16").Value = 66.39189438796498
$ws.Range("S16").Value = 0.004112349560948179
$ws.Range("T16").Value = 0.002799643866861897
$ws.Range("G17").Value = 0.010485
$ws.Range("H17").Value = 0.031455
$ws.Range("I17").Value = 0.001624843498075486
$ws.Range("J17").Value = 0.001634547755634201
$ws.Range("M17").Value = 89.38217433333334
$ws.Range("N17").Value = 268.146523
$ws.Range("O17").Value = 0.2143552015363441
$ws.Range("P17").Value = 0.2175965347165783
$ws.Range("Q17").Value = 0.937172097885
$ws.Range("R17").Value = 8.434548880965
$ws.Range("S17").Value = 0.0003482936554949892
$ws.Range("T17").Value = 0.0003556719274547626
$ws.Range("G18").Value = 0.010485
$ws.Range("H18").Value = 0.031455
$ws.Range("I18").Value = 0.001624843498075486
$ws.Range("J18").Value = 0.001634547755634201
$ws.Range("O18").Value = 0.2934277926151677
$ws.Range("P18").Value = 0.2978648075949286
$ws.Range("Q18").Value = 1.282881581655
$ws.Range("R18").Value = 11.545934234895
$ws.Range("S18").Value = 0.0004767742409853974
$ws.Range("T18").Value = 0.0004868742527367036
$ws.Range("G19").Value = 0.010485
$ws.Range("H19").Value = 0.031455
$ws.Range("I19").Value = 0.001624843498075486
$ws.Range("J19").Value = 0.001634547755634201
$ws.Range("M19").Value = 90.33462533333334
$ws.Range("N19").Value = 271.003876
$ws.Range("O19").Value = 0.2166393574945233
$ws.Range("P19").Value = 0.2199152301234996
$ws.Range("Q19").Value = 0.94715854662
$ws.Range("R19").Value = 8.524426919579998
$ws.Range("S19").Value = 0.000352005051452227
$ws.Range("T19").Value = 0.0003594619458281451
$ws.Range("G20").Value = 0.010485
$ws.Range("H20").Value = 0.031455
$ws.Range("I20").Value = 0.001624843498075486
$ws.Range("J20").Value = 0.001634547755634201
$ws.Range("M20").Value = 18.634161
$ws.Range("N20").Value = 37.268322
$ws.Range("O20").Value = 0.0446882095496985
$ws.Range("P20").Value = 0.03024263611988591
$ws.Range("Q20").Value = 0.195379178085
$ws.Range("R20").Value = 1.17227506851
$ws.Range("S20").Value = 0.00007261134672746247
$ws.Range("T20").Value = 0.00004943303299422134
$ws.Range("G21").Value = 0.010485
$ws.Range("H21").Value = 0.031455
$ws.Range("I21").Value = 0.001624843498075486
$ws.Range("J21").Value = 0.001634547755634201
$ws.Range("M21").Value = 96.27664699999998
$ws.Range("N21").Value = 288.829941
$ws.Range("O21").Value = 0.2308894388042666
$ws.Range("P21").Value = 0.2343807914451077
$ws.Range("Q21").Value = 1.009460643795
$ws.Range("R21").Value = 9.085145794154998
$ws.Range("S21").Value = 0.0003751592034154105
$ws.Range("T21").Value = 0.0003831065966203684
